$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove Sheet2 (user creation sheet no longer needed)
$wb.Worksheets.Item("Sheet2").Delete()

$ws = $wb.Worksheets.Item("Sheet1")

# Fill column A with the full loes1..loes20 user list
for ($i = 1; $i -le 20; $i++) {
    $ws.Cells.Item($i, 1).Value = "loes$i"
}

# Column B gets an (empty, for now) numeric id column next to each user
$ws.Range("B1:B20").NumberFormat = "0"

$ws.Columns("B:C").AutoFit()

$ws.Range("B10").Select()

$wb.Save()
